# Update results with new values (three_way_quicksort results sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the BEST_CASE / AVERAGE_CASE row labels -------------------------
# Row 2 was BEST_CASE, now holds the AVERAGE_CASE series; row 3 was
# AVERAGE_CASE, now holds the BEST_CASE series. WORST_CASE (row 4) stays put.
$ws.Range("A2").Value = "AVERAGE_CASE"
$ws.Range("A3").Value = "BEST_CASE"

# --- New measured values ---------------------------------------------------
$row2 = @(5900, 9900, 17640, 32910, 61960, 123790, 233340, 463230, 916310)
$row3 = @(590, 1100, 1750, 3350, 6980, 13180, 26620, 52880, 106410)
$row4 = @(19270, 73190, 280050, 1104630, 4378480, 17444460, 69603660, 278857410, 1111820500)

for ($i = 0; $i -lt 9; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $row2[$i]
    $ws.Cells.Item(3, 2 + $i).Value = $row3[$i]
    $ws.Cells.Item(4, 2 + $i).Value = $row4[$i]
}

# --- Chart 2 no longer excludes the worst case, simplify its title --------
$co2 = $ws.ChartObjects(2)
$co2.Chart.ChartTitle.Text = "Three-Way-Quicksort"

# --- Restore the cursor / selection left behind by the author -------------
$ws.Range("T21").Select() | Out-Null
